$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114 (shifts existing rows 114-142 down to 115-143),
# then populate the new row with this week's record (same template as the
# previous row 114 record, but with an updated date and origin).
$ws.Rows.Item(114).Insert()

$row = 114
$ws.Cells.Item($row, 1).Value2 = 11
$ws.Cells.Item($row, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value2 = "Bíobío"
$ws.Cells.Item($row, 4).Value2 = 44889
$ws.Cells.Item($row, 5).Value2 = 8
$ws.Cells.Item($row, 6).Value2 = "Fruta"
$ws.Cells.Item($row, 7).Value2 = 100108
$ws.Cells.Item($row, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value2 = 100108002
$ws.Cells.Item($row, 10).Value2 = "Mango"
$ws.Cells.Item($row, 11).Value2 = "Sin especificar"
$ws.Cells.Item($row, 12).Value2 = "Primera"
$ws.Cells.Item($row, 13).Value2 = 200
$ws.Cells.Item($row, 14).Value2 = 8000
$ws.Cells.Item($row, 15).Value2 = 8500
$ws.Cells.Item($row, 16).Value2 = 8250
$ws.Cells.Item($row, 17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item($row, 18).Value2 = "Ecuador"
$ws.Cells.Item($row, 19).Value2 = 2062
$ws.Cells.Item($row, 20).Value2 = 4
